$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "224×5=1120" "318×6=1908"
Replace-Text "823×9=7407" "986×6=5916"
Replace-Text "389×5=1945" "316×9=2844"
Replace-Text "548×6=3288" "922×6=5532"
Replace-Text "512×4=2048" "337×2=674"
Replace-Text "892×4=3568" "539×3=1617"
Replace-Text "605×9=5445" "528×5=2640"
Replace-Text "251×2=502" "616×2=1232"
Replace-Text "523×6=3138" "557×4=2228"
Replace-Text "930×5=4650" "625×9=5625"
Replace-Text "591×4=2364" "155×8=1240"
Replace-Text "662×9=5958" "287×6=1722"
Replace-Text "536×6=3216" "427×5=2135"
Replace-Text "156×7=1092" "207×6=1242"
Replace-Text "123×4=492" "125×5=625"
Replace-Text "303×6=1818" "532×3=1596"
Replace-Text "396×6=2376" "309×2=618"
Replace-Text "927×9=8343" "867×2=1734"
Replace-Text "221×8=1768" "858×2=1716"
Replace-Text "145×5=725" "944×7=6608"
Replace-Text "928×7=6496" "518×8=4144"
Replace-Text "549×7=3843" "522×7=3654"
Replace-Text "707×6=4242" "421×9=3789"
Replace-Text "142×8=1136" "316×6=1896"
Replace-Text "708×4=2832" "991×3=2973"
